$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook is a weekly-refreshed price series. Two new weekly records
# need to be spliced into the existing data block, pushing the subsequent
# rows down (matching the target diff exactly):
#   - a new row inserted at row 252
#   - a new row inserted at row 275 (post first insertion numbering)

function Set-RabanitoRow {
    param([int]$Row, [double]$D, [double]$J, [double]$K, [double]$L, [double]$M, [string]$O, [double]$P)

    $ws.Cells.Item($Row, 1).Value = 9
    $ws.Cells.Item($Row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($Row, 3).Value = "Metropolitana"
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = 13
    $ws.Cells.Item($Row, 6).Value = 300000001
    $ws.Cells.Item($Row, 7).Value = "Rabanito"
    $ws.Cells.Item($Row, 8).Value = "Sin especificar"
    $ws.Cells.Item($Row, 9).Value = "Primera"
    $ws.Cells.Item($Row, 10).Value = $J
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = "`$/cien unidades (volumen en unidades)"
    $ws.Cells.Item($Row, 15).Value = $O
    $ws.Cells.Item($Row, 16).Value = $P
    $ws.Cells.Item($Row, 17).Value = 100
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

# Insert first new record at row 252 (shifts old rows 252..285 down to 253..286)
$ws.Rows.Item(252).EntireRow.Insert()
Set-RabanitoRow 252 44748 7000 2500 3000 2750 "Provincia de Chacabuco" 28

# Insert second new record at row 275 (shifts old rows 274..286 down to 276..287)
$ws.Rows.Item(275).EntireRow.Insert()
Set-RabanitoRow 275 44747 7000 3000 3000 3000 "Provincia de Chacabuco" 30
